$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.786.97'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').Value = '3.312.21'
$ws.Range('E3').Value = '  +6.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '601.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.73'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.78%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.309.15'
$ws.Range('E8').Value = '  +6.40%  '
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('E10').Value = '  +2.84%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.10%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.473'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.09%  '
$ws.Range('E13').Value = '  +1.07%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.80'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').Value = '3.858.36'
$ws.Range('E15').Value = '  +6.29%  '
$ws.Range('D17').Value = '3.310.63'
$ws.Range('E17').Value = '  +6.12%  '
$ws.Range('D18').Value = '63.927.70'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.91'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.84%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '480.99'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('E21').Value = '  +0.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.736'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.01'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +5.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.76'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.71%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.29'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.86%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  +3.20%  '
$ws.Range('E31').Value = '  +4.66%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '29.35'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +9.77%  '
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.11'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '53.03'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0753'
$ws.Range('E38').Value = '  +7.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0404'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.85%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '430.31'
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Value = '3.049.59'
$ws.Range('E41').Value = '  +5.19%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.43'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.90%  '
$ws.Range('E43').Value = '  +2.60%  '
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.267'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '26.48'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +4.14%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '36.18'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +15.30%  '
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('E51').Value = '  +2.46%  '
